$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 currently holds "Luiz Felipe" stats, row 21 holds "Florian Lejeune" stats.
# The edit swaps these two players' stat rows (columns C through DJ), while
# columns A, B (League/Team) and DK, DL (tag/extra) stay untouched since they
# are identical for both rows already.

$firstCol = 3   # column C
$lastCol  = 114 # column DJ

$row20Range = $ws.Range($ws.Cells.Item(20, $firstCol), $ws.Cells.Item(20, $lastCol))
$row21Range = $ws.Range($ws.Cells.Item(21, $firstCol), $ws.Cells.Item(21, $lastCol))

$row20Values = $row20Range.Value2
$row21Values = $row21Range.Value2

$row20Range.Value2 = $row21Values
$row21Range.Value2 = $row20Values
